# daily auto push: 2026-01-14 03:53 UTC
# Insert a new timestamp row for 2026/01/14 03:53 UTC (hour bucket 11)
# just above the existing 2026/12/29 block, shifting every following row
# down by one (old row 635 -> 636, ... old row 676 -> 677).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 634 already holds the same date/weekday text ("2026/01/14", "水");
# copy it so the new row inherits the identical cell typing/formatting,
# then insert the copy above row 635 (pushing the 12/29 block down).
$ws.Rows.Item(634).Copy()
$ws.Rows.Item(635).Insert()

# Only the "time bucket" value differs for the new row.
$ws.Cells.Item(635, 3).Value = 11
